$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in cell A1 by one day (serial 45310 -> 45311,
# i.e. 2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# Update the two prices that changed (fix for the exceeded Google Drive
# request bug mentioned in the commit message)
$ws.Range("D30").Value = 516
$ws.Range("D31").Value = 999
